$d = $word.ActiveDocument

# Insert a new Title-styled paragraph ("  040 item list") as the very first
# paragraph of the document body, ahead of the existing list items.
$target = $d.Range(0, 0)

$titleParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve"> 040 item list</w:t></w:r></w:p>'

[void]$target.InsertXML($titleParagraphXml)

Write-Host "Inserted title paragraph; document now has $($d.Paragraphs.Count) paragraphs."
